$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column before the "Late" column (N) to make room for a new
# (for now blank) field - part of the RBI / Variable Instalments layout
# change. This shifts "Late" (N->O) and "Outstanding" (P->Q) one column right.
$ws.Columns.Item(14).Insert()

# The newly inserted column should take on the same width as its neighbour
# column M ("In Advance"), matching the final layout.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Update the active selection as recorded after the edit.
$ws.Range("T8").Select()
